# The "Audio" column (column A, containing the header "Audio" and the
# merged "audio.mp3" cell spanning A2:A4) is removed entirely. Every other
# column shifts one place to the left (old B -> A, old C -> B, ... old G -> F)
# along with its formatting, the A2:A4 merge disappears (it belonged to the
# deleted column), and the column-G width=16 setting now belongs to column F.
#
# This is exactly what Excel does when a user selects the whole of column A
# and deletes it, so that's what we drive via COM.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select column A (mirrors the author selecting the whole column before
# deleting it) and then delete it, shifting everything else left.
$ws.Range("A:A").Select()
$ws.Range("A:A").EntireColumn.Delete()
